$wb = $excel.ActiveWorkbook

# Work on the "RequestBuyersList" worksheet (sheet5.xml) - add a new "Company" column
$ws = $wb.Worksheets.Item("RequestBuyersList")

# Make RequestBuyersList the active (selected) tab of the workbook
$ws.Activate()

$ws.Range("D1").Value = "Company"
$ws.Range("D2").Value = "StandardTestCompany"

# Match the header formatting used by the other header cells (A1:C1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection on this sheet to E2
$ws.Range("E2").Select()

$wb.Save()
